$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the existing row 1111.
# This shifts the old rows 1111-1212 down to 1113-1214, which already
# reproduces the tail of the target sheet (new rows 1213/1214) with no
# further work needed.
$ws.Rows("1111:1112").Insert()

# Populate new row 1111 (newest "1a (guarda)" price entry, dated 2022-08-10)
$ws.Cells.Item(1111, 1).Value = 9
$ws.Cells.Item(1111, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1111, 3).Value = "Metropolitana"
$ws.Cells.Item(1111, 4).Value = 44783
$ws.Cells.Item(1111, 5).Value = 13
$ws.Cells.Item(1111, 6).Value = 100112045
$ws.Cells.Item(1111, 7).Value = "Zapallo"
$ws.Cells.Item(1111, 8).Value = "Camote"
$ws.Cells.Item(1111, 9).Value = "1a (guarda)"
$ws.Cells.Item(1111, 10).Value = 180
$ws.Cells.Item(1111, 11).Value = 1000
$ws.Cells.Item(1111, 12).Value = 1100
$ws.Cells.Item(1111, 13).Value = 1050
$ws.Cells.Item(1111, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(1111, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(1111, 16).Value = 1050
$ws.Cells.Item(1111, 17).Value = 1
$ws.Cells.Item(1111, 18).Value = "Hortaliza"

# Populate new row 1112 (newest "2a (guarda)" price entry, dated 2022-08-10)
$ws.Cells.Item(1112, 1).Value = 9
$ws.Cells.Item(1112, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1112, 3).Value = "Metropolitana"
$ws.Cells.Item(1112, 4).Value = 44783
$ws.Cells.Item(1112, 5).Value = 13
$ws.Cells.Item(1112, 6).Value = 100112045
$ws.Cells.Item(1112, 7).Value = "Zapallo"
$ws.Cells.Item(1112, 8).Value = "Camote"
$ws.Cells.Item(1112, 9).Value = "2a (guarda)"
$ws.Cells.Item(1112, 10).Value = 70
$ws.Cells.Item(1112, 11).Value = 780
$ws.Cells.Item(1112, 12).Value = 800
$ws.Cells.Item(1112, 13).Value = 790
$ws.Cells.Item(1112, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(1112, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(1112, 16).Value = 790
$ws.Cells.Item(1112, 17).Value = 1
$ws.Cells.Item(1112, 18).Value = "Hortaliza"
